# The deck currently uses the "Integral" design theme (ppt/theme/theme2.xml,
# wired up via the slide master / presentation relationships) while
# ppt/theme/theme1.xml (only used by the notes master) still holds the
# original default "Office Theme" palette.
#
# The target edit swaps the two themes' content: the live/visible design
# theme becomes the classic "Office Theme" color palette, while the
# notes-only theme keeps/holds the old "Integral" colors.
#
# The notes-master theme part isn't reachable/writable through the
# PowerPoint COM object model (there is no supported automation surface to
# author a distinct notes-master theme), so we drive the change through the
# one theme objects that IS reachable and authoritative for the design
# theme actually used by the slides: Slide.ThemeColorScheme.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index order (matches PowerPoint's ppThemeColorXXX / ThemeColorScheme.Item):
#   1 dk1   2 lt1   3 dk2   4 lt2
#   5 accent1  6 accent2  7 accent3  8 accent4  9 accent5  10 accent6
#   11 hlink   12 folHlink
#
# RGB values are packed as VBA's RGB(r,g,b) = r + g*256 + b*65536 (i.e. the
# COM RGBColor.RGB property is BGR-ordered), matching the Office Theme's
# default color scheme:
#   dk1=000000 lt1=FFFFFF dk2=44546A lt2=E7E6E6
#   accent1=5B9BD5 accent2=ED7D31 accent3=A5A5A5 accent4=FFC000
#   accent5=4472C4 accent6=70AD47 hlink=0563C1 folHlink=954F72

$officeThemeRgb = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeThemeRgb[$i - 1]
}
